$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.05
$ws.Range("I2").Value = 3.25
$ws.Range("L2").Value = 3.6
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.53
$ws.Range("W2").Value = 1.3
$ws.Range("X2").Value = 3.4
$ws.Range("AE2").Value = 17
$ws.Range("G3").Value = 2.5
$ws.Range("I3").Value = 2.7
$ws.Range("J3").Value = 3.2
$ws.Range("L3").Value = 3.25
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 2.02
$ws.Range("AB3").Value = 15
$ws.Range("AD3").Value = 29
$ws.Range("AM3").Value = 15
$ws.Range("AN3").Value = 11
$ws.Range("AO3").Value = 29
$ws.Range("G4").Value = 2.6
$ws.Range("I4").Value = 2.63
$ws.Range("G6").Value = 1.98
$ws.Range("H6").Value = 3.55
$ws.Range("I6").Value = 3.45
$ws.Range("J6").Value = 2.57
$ws.Range("K6").Value = 2.18
$ws.Range("L6").Value = 3.85
$ws.Range("P6").Value = 3.55
$ws.Range("Z6").Value = 2.1
$ws.Range("AA6").Value = 8
$ws.Range("AB6").Value = 9.75
$ws.Range("AD6").Value = 17.5
$ws.Range("AE6").Value = 15
$ws.Range("AH6").Value = 6.8
$ws.Range("AI6").Value = 13.5
$ws.Range("AJ6").Value = 55
$ws.Range("AK6").Value = 400
$ws.Range("AL6").Value = 11.75
$ws.Range("AM6").Value = 19.5
$ws.Range("AN6").Value = 11.75
$ws.Range("AO6").Value = 50
$ws.Range("AP6").Value = 28
$ws.Range("AQ6").Value = 32
$ws.Range("G10").Value = 1.62
$ws.Range("J10").Value = 2.25
$ws.Range("K10").Value = 2.2
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 10
$ws.Range("AK10").Value = 451
$ws.Range("Q11").Value = 1.95
$ws.Range("R11").Value = 1.9
$ws.Range("G12").Value = 2.9
$ws.Range("AC12").Value = 12
$ws.Range("AD12").Value = 34
$ws.Range("S13").Value = 2.03
$ws.Range("T13").Value = 1.78
$ws.Range("AK13").Value = 126
$ws.Range("G15").Value = 4.2
$ws.Range("I15").Value = 1.7
$ws.Range("M15").Value = 1.04
$ws.Range("N15").Value = 9
$ws.Range("AM15").Value = 8.5
$ws.Range("L16").Value = 10
$ws.Range("N16").Value = 15
$ws.Range("O16").Value = 1.08
$ws.Range("P16").Value = 7
$ws.Range("U16").Value = 1.73
$ws.Range("V16").Value = 2
$ws.Range("Y16").Value = 1.91
$ws.Range("Z16").Value = 1.8
$ws.Range("AC16").Value = 12
$ws.Range("AJ16").Value = 67
$ws.Range("AP16").Value = 81
$ws.Range("M17").Value = 1.02
$ws.Range("N17").Value = 11
$ws.Range("M18").Value = 1.05
$ws.Range("N18").Value = 8.5
$ws.Range("Q18").Value = 1.83
$ws.Range("R18").Value = 1.98
$ws.Range("AG18").Value = 11
$ws.Range("AO18").Value = 29
$ws.Range("G23").Value = 3.8
$ws.Range("H23").Value = 3.75
$ws.Range("I23").Value = 1.85
$ws.Range("J23").Value = 4.1
$ws.Range("L23").Value = 2.4
$ws.Range("M23").Value = 1.05
$ws.Range("N23").Value = 8.75
$ws.Range("O23").Value = 1.24
$ws.Range("P23").Value = 3.75
$ws.Range("Q23").Value = 1.72
$ws.Range("R23").Value = 2.05
$ws.Range("U23").Value = 2.75
$ws.Range("V23").Value = 1.42
$ws.Range("Y23").Value = 1.65
$ws.Range("Z23").Value = 2.1
$ws.Range("AA23").Value = 11.75
$ws.Range("AB23").Value = 23
$ws.Range("AC23").Value = 13.5
$ws.Range("AD23").Value = 60
$ws.Range("AE23").Value = 35
$ws.Range("AF23").Value = 40
$ws.Range("AG23").Value = 8.75
$ws.Range("AI23").Value = 15
$ws.Range("AJ23").Value = 60
$ws.Range("AK23").Value = 450
$ws.Range("AL23").Value = 8
$ws.Range("AM23").Value = 10
$ws.Range("AP23").Value = 15
$ws.Range("AQ23").Value = 25
